$wb = $excel.ActiveWorkbook

# --- "Test Cases" sheet: add a "Result" column (D) ---
$wsCases = $wb.Worksheets.Item("Test Cases")
$wsCases.Range("D1").Value = "Result"
$wsCases.Range("D1").Font.Bold = $true
$wsCases.Range("D2").Select()

# --- "Test Steps" sheet: add a "Result" column (F) ---
$wsSteps = $wb.Worksheets.Item("Test Steps")
$wsSteps.Range("F1").Value = "Result"
$wsSteps.Range("F1").Font.Bold = $true
$wsSteps.Range("F19").Select()

# "Test Steps" becomes the active/visible tab
$wsSteps.Activate()
